$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '57.840.61'
Set-TextValue 2 5 '  -1.72%  '

Set-TextValue 3 4 '2.453.39'
Set-TextValue 3 5 '  -1.70%  '

Set-TextValue 4 5 '  -0.24%  '

Set-TextValue 5 4 '517.81'
Set-TextValue 5 5 '  -3.44%  '

Set-TextValue 6 4 '132.09'
Set-TextValue 6 5 '  -2.97%  '

Set-TextValue 7 4 '1.00'
Set-TextValue 7 5 '  +0.09%  '

Set-TextValue 8 5 '  -1.92%  '

Set-TextValue 9 4 '2.457.62'
Set-TextValue 9 5 '  -2.53%  '

Set-TextValue 10 5 '  -3.32%  '

Set-TextValue 11 5 '  -0.31%  '

Set-TextValue 12 4 '5.26'
Set-TextValue 12 5 '  -1.21%  '

Set-TextValue 13 5 '  -2.69%  '

Set-TextValue 14 4 '2.890.96'
Set-TextValue 14 5 '  -2.31%  '

Set-TextValue 15 4 '57.783.99'
Set-TextValue 15 5 '  -1.82%  '

Set-TextValue 16 4 '22.19'
Set-TextValue 16 5 '  -3.38%  '

Set-TextValue 17 4 '0.0000135'
Set-TextValue 17 5 '  -2.80%  '

Set-TextValue 18 4 '2.457.33'
Set-TextValue 18 5 '  -2.64%  '

Set-TextValue 19 4 '10.61'
Set-TextValue 19 5 '  -3.94%  '

Set-TextValue 20 4 '318.80'
Set-TextValue 20 5 '  -1.21%  '

Set-TextValue 21 4 '4.15'
Set-TextValue 21 5 '  -2.56%  '

Set-TextValue 22 5 '  -0.05%  '

Set-TextValue 23 5 '  -4.34%  '

Set-TextValue 24 4 '64.30'
Set-TextValue 24 5 '  -0.91%  '

Set-TextValue 25 4 '0.408'
Set-TextValue 25 5 '  -2.78%  '

Set-TextValue 26 4 '0.999'
Set-TextValue 26 5 '  +0.04%  '

Set-TextValue 27 5 '  -2.72%  '

Set-TextValue 28 4 '7.34'
Set-TextValue 28 5 '  -2.38%  '

Set-TextValue 29 4 '0.0₃0739'
Set-TextValue 29 5 '  -3.49%  '

Set-TextValue 30 4 '167.35'
Set-TextValue 30 5 '  -1.90%  '

Set-TextValue 31 5 '  -4.21%  '

Set-TextValue 32 4 '6.20'
Set-TextValue 32 5 '  -6.57%  '

Set-TextValue 33 4 '1.16'
Set-TextValue 33 5 '  -0.76%  '

Set-TextValue 34 5 '  -0.01%  '

Set-TextValue 36 4 '18.04'
Set-TextValue 36 5 '  -1.67%  '

Set-TextValue 37 4 '1.30'
Set-TextValue 37 5 '  -5.53%  '

Set-TextValue 38 4 '3.96'
Set-TextValue 38 5 '  -2.46%  '

Set-TextValue 39 5 '  -4.20%  '

Set-TextValue 40 4 '36.21'
Set-TextValue 40 5 '  -1.89%  '

Set-TextValue 41 4 '0.785'
Set-TextValue 41 5 '  -3.02%  '

Set-TextValue 42 4 '3.43'
Set-TextValue 42 5 '  -4.16%  '

Set-TextValue 43 4 '270.45'
Set-TextValue 43 5 '  -4.77%  '

Set-TextValue 44 4 '4.94'
Set-TextValue 44 5 '  -4.29%  '

Set-TextValue 45 4 '0.585'
Set-TextValue 45 5 '  -3.40%  '

Set-TextValue 46 4 '124.86'
Set-TextValue 46 5 '  -4.03%  '

Set-TextValue 47 4 '0.0908'
Set-TextValue 47 5 '  -1.61%  '

Set-TextValue 48 4 '0.0485'
Set-TextValue 48 5 '  -3.53%  '

Set-TextValue 49 4 '0.0211'
Set-TextValue 49 5 '  -3.75%  '

Set-TextValue 50 4 '16.74'
Set-TextValue 50 5 '  -3.56%  '

Set-TextValue 51 4 '1.724.91'
Set-TextValue 51 5 '  -1.82%  '
